# Apply updates to the "想去人数" (number of people interested) counts
# across the "展览", "演出" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 584
$ws1.Range("F6").Value = 349
$ws1.Range("F7").Value = 1618

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 103

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 584
$ws4.Range("F6").Value = 349
$ws4.Range("F7").Value = 103
$ws4.Range("F11").Value = 1618
